$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 16: fill in previously-empty hours cells and adjust downstream inputs
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 2
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1
$ws.Range("Q16").Value = 15
$ws.Range("S16").Value = 100
$ws.Range("T16").ClearContents()
$ws.Range("X16").Formula = '=SUM($W$3:W16)+SUM(V17:$V$23)'

# Row 18: T18 now carries the 100 that used to live on row 16
$ws.Range("T18").Value = 100

# Move the active selection / visible top-left cell to match the saved view
$ws.Range("X16").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 8
